$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A88").Value = "2025-04-29 14:47:42"
$ws.Range("B88").Value = 267
